$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Anotacoes")
$ws2 = $wb.Worksheets.Item("Tabelas")

# --- New header row (row 8): Model / Controller / View in columns C, D, E ---
$ws1.Range("C8").Value = "Model"
$ws1.Range("D8").Value = "Controller"
$ws1.Range("E8").Value = "View"

# Reuse the existing "top box" border style (left+right+top, no bottom) already
# used by the Tabelas sheet's table headers (e.g. B2), by copying its format.
$ws2.Range("B2").Copy()
$ws1.Range("C8:E8").PasteSpecial(-4122)   # xlPasteFormats

# --- Rows 9 and 10: box every cell from B to E with a thin border on all sides ---
# Create the new "plain full box" style once (fontId 0 + full thin border), then
# propagate it by copying the format so only a single new style gets minted.
$ws1.Range("B9").Borders.LineStyle = 1
$ws1.Range("B9").Copy()
$ws1.Range("C9:E9").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("B10:E10").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# --- Column widths for the new C/D/E columns ---
$ws1.Columns.Item(3).ColumnWidth = 5.8
$ws1.Columns.Item(4).ColumnWidth = 9.1
$ws1.Columns.Item(5).ColumnWidth = 6.26

# --- Update active selection to C10, matching the post-edit cursor position ---
$ws1.Range("C10").Select()
